# Fixed #418 Empty AQL expressions generate empty lines.
#
# The document contains an empty paragraph (just an empty run <w:t/>)
# right after the "Start of demonstration:" paragraph. This empty
# paragraph is removed entirely so the "Some value" paragraph directly
# follows "Start of demonstration:".

$d = $word.ActiveDocument

# Locate the empty paragraph: it is the second paragraph in the body
# and its range text is empty (only the paragraph mark).
$target = $null
foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    # Paragraph mark is represented by a single trailing char; an empty
    # paragraph's Range.Text is just that mark (length 1) or empty.
    $trimmed = $text.TrimEnd([char]13, [char]7)
    if ($trimmed.Length -eq 0) {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.Delete()
}
